$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item('展览')
$ws.Range('F2').Value = 292
$ws.Range('F3').Value = 996
$ws.Range('F4').Value = 1232
$ws.Range('F6').Value = 3227
$ws.Range('I9').Value = '//i0.hdslb.com/bfs/openplatform/202404/t2hwbRC01713235949385.jpeg'
$ws.Range('F10').Value = 705
$ws.Range('F11').Value = 560
$ws.Range('F13').Value = 36
$ws.Range('F14').Value = 81
$ws.Range('F15').Value = 637
$ws.Range('F16').Value = 1437
$ws.Range('F17').Value = 1437
$ws.Range('F20').Value = 583
$ws.Range('F21').Value = 290
$ws.Range('F23').Value = 485
$ws.Range('F24').Value = 29810
$ws.Range('F25').Value = 29820
$ws.Range('F26').Value = 698
$ws.Range('F27').Value = 628
$ws.Range('F28').Value = 23096
$ws.Range('F29').Value = 23143
$ws.Range('F30').Value = 382
$ws.Range('F32').Value = 35
$ws.Range('F34').Value = 196
$ws.Range('F36').Value = 452
$ws.Range('F37').Value = 1123
$ws.Range('F38').Value = 5248
$ws.Range('F39').Value = 657
$ws.Range('F40').Value = 388
$ws.Range('F42').Value = 297

# ---- Sheet: 演出 ----
$ws = $wb.Worksheets.Item('演出')
$ws.Range('F5').Value = 209
$ws.Range('F18').Value = 4
$ws.Range('C21').Value = '上海·爵士情调女王KAREN SOUZA凯伦索萨2024演唱会'
$ws.Range('D21').Value = '南京西路1376号上海商城4层 商城剧院'
$ws.Range('E21').Value = '2024.05.03 19:30-05.03 21:00'
$ws.Range('F21').Value = 4
$ws.Range('G21').Value = 280
$ws.Range('H21').Value = 'https://show.bilibili.com/platform/detail.html?id=82653'
$ws.Range('I21').Value = '//i0.hdslb.com/bfs/openplatform/202403/IkBVehui1710141982443.jpeg'
$ws.Range('C22').Value = '上海·申放送-Virtual Super Live-2024 in Shanghai'
$ws.Range('D22').Value = '中兴路1599号金融街融泰中心 蜚声上海PHASE LIVE HOUSE'
$ws.Range('E22').Value = '2024.05.03 19:00-05.03 22:00'
$ws.Range('F22').Value = 446
$ws.Range('G22').Value = 488
$ws.Range('H22').Value = 'https://show.bilibili.com/platform/detail.html?id=83102'
$ws.Range('I22').Value = '//i0.hdslb.com/bfs/openplatform/202404/lhju6VbJ1712475891713.jpeg'
$ws.Range('F34').Value = 954
$ws.Range('F35').Value = 478
$ws.Range('F37').Value = 63
$ws.Range('F38').Value = 63
$ws.Range('F45').Value = 21
$ws.Range('F47').Value = 3

# ---- Sheet: 本地生活 ----
$ws = $wb.Worksheets.Item('本地生活')
$ws.Range('F4').Value = 683

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item('全部类型')
$ws.Range('F2').Value = 683
$ws.Range('F3').Value = 292
$ws.Range('F6').Value = 996
$ws.Range('F7').Value = 1232
$ws.Range('F8').Value = 209
$ws.Range('F10').Value = 3227
$ws.Range('I13').Value = '//i0.hdslb.com/bfs/openplatform/202404/t2hwbRC01713235949385.jpeg'
$ws.Range('F14').Value = 705
$ws.Range('F18').Value = 560
$ws.Range('F20').Value = 36
$ws.Range('F21').Value = 81
$ws.Range('F22').Value = 637
$ws.Range('F23').Value = 1437
$ws.Range('F24').Value = 1437
$ws.Range('C27').Value = '上海·百梦动漫游戏嘉年华'
$ws.Range('D27').Value = '吴中路1588号上海爱琴海购物中心F4 百忍潮玩对战中心(爱琴海购物中心店)'
$ws.Range('E27').Value = '2024.05.01 10:00-05.02 19:00'
$ws.Range('F27').Value = 8
$ws.Range('G27').Value = 49
$ws.Range('H27').Value = 'https://show.bilibili.com/platform/detail.html?id=84152'
$ws.Range('I27').Value = '//i2.hdslb.com/bfs/openplatform/202404/TJknSP7V1712849614164.jpeg'
$ws.Range('C28').Value = '上海·第五十七届燃梦星辰动漫嘉年华'
$ws.Range('D28').Value = '云锦路500号(近11号线地铁站5号口) 绿地滨江CLUB'
$ws.Range('E28').Value = '2024.05.01 10:30-05.01 16:30'
$ws.Range('F28').Value = 583
$ws.Range('G28').Value = 58.8
$ws.Range('H28').Value = 'https://show.bilibili.com/platform/detail.html?id=83807'
$ws.Range('I28').Value = '//i2.hdslb.com/bfs/openplatform/202404/RGLpPX211712156496032.jpeg'
$ws.Range('C29').Value = '上海·街舞音乐剧《时光代理人：法则游戏》'
$ws.Range('D29').Value = '牛庄路704号 中国大戏院'
$ws.Range('E29').Value = '2024.05.01 19:30-05.19 21:00'
$ws.Range('F29').Value = 384
$ws.Range('G29').Value = 188
$ws.Range('H29').Value = 'https://show.bilibili.com/platform/detail.html?id=82995'
$ws.Range('I29').Value = '//i1.hdslb.com/bfs/openplatform/202403/p9ZC2azX1710816437198.png'
$ws.Range('B30').Value = '''2024-05-01'
$ws.Range('B30').Style = 'Normal'
$ws.Range('C30').Value = '上海·魔都劳动节漫展-CF01'
$ws.Range('D30').Value = '澳门路168号 月星家居（澳门路）'
$ws.Range('E30').Value = '2024.05.01 10:00-05.05 16:00'
$ws.Range('F30').Value = 290
$ws.Range('G30').Value = 49
$ws.Range('H30').Value = 'https://show.bilibili.com/platform/detail.html?id=82992'
$ws.Range('I30').Value = '//i2.hdslb.com/bfs/openplatform/202403/I7O9LMtb1710752670542.jpeg'
$ws.Range('C31').Value = '上海·2024GAF插画艺术节'
$ws.Range('D31').Value = '博成路850号 上海世博展览馆'
$ws.Range('E31').Value = '2024.05.02 10:30-05.04 19:00'
$ws.Range('F31').Value = 485
$ws.Range('G31').Value = 128
$ws.Range('H31').Value = 'https://show.bilibili.com/platform/detail.html?id=83699'
$ws.Range('I31').Value = '//i1.hdslb.com/bfs/openplatform/202403/APlNld8y1711825700811.jpeg'
$ws.Range('F32').Value = 29825
$ws.Range('F34').Value = 698
$ws.Range('F35').Value = 628
$ws.Range('F36').Value = 23245
$ws.Range('F37').Value = 382
$ws.Range('F39').Value = 196
$ws.Range('F41').Value = 452
$ws.Range('F42').Value = 1123
$ws.Range('F43').Value = 5248
$ws.Range('F44').Value = 657
$ws.Range('F45').Value = 478
$ws.Range('F46').Value = 388
$ws.Range('F47').Value = 63
$ws.Range('F49').Value = 297
$ws.Range('F54').Value = 3

